$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 15, shifting existing rows 15-43 down to 16-44.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(15, 3).Value = "Bíobío"
$ws.Cells.Item(15, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(15, 4).Style = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Cells.Item(15, 5).Value = 8
$ws.Cells.Item(15, 6).Value = 100112012
$ws.Cells.Item(15, 7).Value = "Espinaca"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 6000
$ws.Cells.Item(15, 12).Value = 6500
$ws.Cells.Item(15, 13).Value = 6250
$ws.Cells.Item(15, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 625
$ws.Cells.Item(15, 17).Value = 10
$ws.Cells.Item(15, 18).Value = "Hortaliza"
